# "pushign the changes doen for the failures in jenkins-Shaheena"
#
# Updates the Registration/Login test-data rows to a new test user
# (DeltaFaucet8) and flips the previously-failing results to PASS, then
# removes the now-unused "Result1" helper column from the Test Steps sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Registration sheet: new test user + passing result
# ---------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Registration")
$wsReg.Range("A2").Value = "'DeltaFaucet8"
$wsReg.Range("B2").Value = "'Test8"
$wsReg.Range("C2").Value = "'DeltaFaucet2@gmail.com"
$wsReg.Range("F2").Value = "'PASS"
$wsReg.Columns("F").ColumnWidth = 29.5703125
$wsReg.Range("D6").Select() | Out-Null

# ---------------------------------------------------------------------
# Login sheet: same new test user + passing verification message
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A2").Value = "'DeltaFaucet8"
$wsLogin.Range("B2").Value = "'DeltaFaucet8@gmail.com"
$wsLogin.Range("E2").Value = "'PASS -- text verified DeltaFaucet7 -- DeltaFaucet7"
$wsLogin.Columns("E").ColumnWidth = 36
$wsLogin.Range("C10").Select() | Out-Null

# ---------------------------------------------------------------------
# Test Steps sheet: drop the now-unused "Result1" column (H)
# ---------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Columns("H").Delete() | Out-Null
$wsSteps.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore original active sheet/selection (Test Cases was untouched)
# ---------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Activate()
$wsCases.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# Best-effort: reflect the enlarged workbook window geometry recorded
# in the saved view state.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 1024.5
$win.Height = 339.75
